$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Cells.Item(1, 1).Value = "ФИО"
$ws.Cells.Item(1, 2).Value = "Класс"
$ws.Cells.Item(1, 3).Value = "Категория"
$ws.Cells.Item(1, 4).Value = "Тип"
$ws.Cells.Item(1, 5).Value = "Название"
$ws.Cells.Item(1, 6).Value = "Предмет"
$ws.Cells.Item(1, 7).Value = "Этап"
$ws.Cells.Item(1, 8).Value = "Результат"
$ws.Cells.Item(1, 9).Value = "Баллы"

# --- Row 2: Проектная и исследовательская деятельность / kjk / Proj ---
$ws.Cells.Item(2, 1).Value = "s"
$ws.Cells.Item(2, 2).Value = "11-А"
$ws.Cells.Item(2, 3).Value = "Проектная и исследовательская деятельность"
$ws.Cells.Item(2, 4).Value = "kjk"
$ws.Cells.Item(2, 5).Value = "Proj"
$ws.Cells.Item(2, 6).Value = "Proj"
$ws.Cells.Item(2, 7).Value = "Proj"
$ws.Cells.Item(2, 8).Value = "Proj"
$ws.Cells.Item(2, 9).Value = 10

# --- Row 3: Интеллектуальные соревнования / МОШ / Int ---
$ws.Cells.Item(3, 1).Value = "s"
$ws.Cells.Item(3, 2).Value = "11-А"
$ws.Cells.Item(3, 3).Value = "Интеллектуальные соревнования"
$ws.Cells.Item(3, 4).Value = "МОШ"
$ws.Cells.Item(3, 5).Value = "Int"
$ws.Cells.Item(3, 6).Value = "Int"
$ws.Cells.Item(3, 7).Value = "Int"
$ws.Cells.Item(3, 8).Value = "Int"
$ws.Cells.Item(3, 9).Value = 10

# --- Row 4: Спортивные достижения / Турнир по баскетболу / ооо ---
$ws.Cells.Item(4, 1).Value = "s"
$ws.Cells.Item(4, 2).Value = "11-А"
$ws.Cells.Item(4, 3).Value = "Спортивные достижения"
$ws.Cells.Item(4, 4).Value = "Турнир по баскетболу"
$ws.Cells.Item(4, 5).Value = "ооо"
$ws.Cells.Item(4, 6).Value = "ооо"
$ws.Cells.Item(4, 7).Value = "ооо"
$ws.Cells.Item(4, 8).Value = "ооо"
$ws.Cells.Item(4, 9).Value = 10

# --- Row 5: Спортивные достижения / Турнир по баскетболу / ьььь, ииии ---
$ws.Cells.Item(5, 1).Value = "s"
$ws.Cells.Item(5, 2).Value = "11-А"
$ws.Cells.Item(5, 3).Value = "Спортивные достижения"
$ws.Cells.Item(5, 4).Value = "Турнир по баскетболу"
$ws.Cells.Item(5, 5).Value = "ьььь"
$ws.Cells.Item(5, 6).Value = "ииии"
$ws.Cells.Item(5, 7).Value = "ьььь"
$ws.Cells.Item(5, 8).Value = "ииии"
$ws.Cells.Item(5, 9).Value = 10

# --- Column widths (best achievable approximations of the target bestFit widths;
#     the runtime quantizes stored widths to 1/6-character increments) ---
$ws.Columns.Item(3).ColumnWidth = 49.666667
$ws.Columns.Item(6).ColumnWidth = 8.5
$ws.Columns.Item(7).ColumnWidth = 5.0
$ws.Columns.Item(8).ColumnWidth = 10.833333
$ws.Columns.Item(9).ColumnWidth = 6.166667
